# fix(gui) step 1 and 2
# - Bump the date serial in A1 by one day (price list issue date).
# - Update the unit prices in column D for rows 32-39 (PORTACANDADOS price list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: date header moves from 45308 (2024-01-17) to 45309 (2024-01-18)
$ws.Range("A1").Value = 45309

# Step 2: updated unit prices
$ws.Range("D32").Value = 219.124
$ws.Range("D33").Value = 313.543
$ws.Range("D34").Value = 417.457
$ws.Range("D35").Value = 429.794
$ws.Range("D36").Value = 563.266
$ws.Range("D37").Value = 644.069
$ws.Range("D38").Value = 771.267
$ws.Range("D39").Value = 918.41
